$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - blog.pabii.co.kr entry
$ws.Range("D9").Value = "통계학 석사가 본 데이터 사이언스"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ms-stat-data-science/#utm_source=rss&utm_medium=rss&utm_campaign=ms-stat-data-science"

# Row 20 - ai-creator.tistory.com entry
$ws.Range("D20").Value = "[AI] 코딩 없이 마스크 착용 감지 딥러닝 서비스 만들기 (w/ Teachable Machine)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/548"

# Row 43 - nittaku.tistory.com entry
$ws.Range("D43").Value = "lg gram에 삼성 노트 PC 사용하기"
$ws.Range("E43").Value = "https://nittaku.tistory.com/510"

# Row 51 - bskyvision.com entry
$ws.Range("D51").Value = "티스토리 회원 분들 필독(댓글, 방명록, 구독신청 전에)"
$ws.Range("E51").Value = "https://bskyvision.com/notice/1178"
